$p = $ppt.ActivePresentation
$s = $p.Slides.Item(8)
$s.Shapes.Item("TextBox 21").Delete()
